$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.591.30"
$ws.Range("E2").Value = "  +2.34%  "

$ws.Range("D3").Value = "2.952.39"
$ws.Range("E3").Value = "  +2.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.36"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.03"
$ws.Range("E6").Value = "  +4.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "2.952.26"
$ws.Range("E8").Value = "  +2.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.506"
$ws.Range("E9").Value = "  +2.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.95"
$ws.Range("E10").Value = "  +0.42%  "

$ws.Range("E11").Value = "  +8.14%  "

$ws.Range("E12").Value = "  +1.29%  "

$ws.Range("E13").Value = "  +6.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.17"
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.125"
$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("D16").Value = "3.443.20"
$ws.Range("E16").Value = "  +2.28%  "

$ws.Range("D17").Value = "62.613.80"
$ws.Range("E17").Value = "  +2.52%  "

$ws.Range("D18").Value = "2.951.62"
$ws.Range("E18").Value = "  +1.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.65"
$ws.Range("E19").Value = "  +1.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "434.40"
$ws.Range("E20").Value = "  +1.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.47"
$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.661"
$ws.Range("E22").Value = "  +0.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.95"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("E24").Value = "  +5.79%  "

$ws.Range("E25").Value = "  -0.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.90"
$ws.Range("E26").Value = "  +4.05%  "

$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("E29").Value = "  +6.13%  "

$ws.Range("E30").Value = "  +3.12%  "

$ws.Range("E31").Value = "  +1.63%  "

$ws.Range("E32").Value = "  +16.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.109"
$ws.Range("E33").Value = "  +2.83%  "

$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("E36").Value = "  +1.19%  "

$ws.Range("E37").Value = "  +1.64%  "

$ws.Range("E38").Value = "  +5.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.67"
$ws.Range("E39").Value = "  +0.70%  "

$ws.Range("E40").Value = "  +3.95%  "

$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").Value = "  -3.21%  "

$ws.Range("E43").Value = "  +3.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.14"
$ws.Range("E44").Value = "  -4.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "134.69"
$ws.Range("E45").Value = "  +1.02%  "

$ws.Range("D46").Value = "2.682.70"
$ws.Range("E46").Value = "  +0.85%  "

$ws.Range("E47").Value = "  +0.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "353.14"
$ws.Range("E48").Value = "  +2.20%  "

$ws.Range("E50").Value = "  +1.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.50"
$ws.Range("E51").Value = "  -1.54%  "
